$wb = $excel.ActiveWorkbook
$gainers = $wb.Worksheets.Item("Top Gainers")
$losers = $wb.Worksheets.Item("Top Losers")

# Full data refresh for "Top Gainers" rows 2-76: columns B (Stock), C (Latest), D (Weekly), E (Monthly)
$gainersData = @(
    @(2, 'SOLARWORLD', 14.7722, 10.7603, 6.2751),
    @(3, 'BLUEDART', 14.4379, 13.8622, 11.117),
    @(4, 'ADANIGREEN', 11.9896, 9.068, 9.5194),
    @(5, 'DREDGECORP', 11.679, 16.1655, 16.9249),
    @(6, 'INFOBEAN', 9.9924, 23.0892, 38.1616),
    @(7, 'VBL', 9.3912, 7.6723, 11.9675),
    @(8, 'HEG', 8.19, 12.197, 14.5969),
    @(9, 'BUTTERFLY', 7.8152, 10.7306, 13.3343),
    @(10, 'FIVESTAR', 7.6902, 7.7002, 7.7805),
    @(11, 'WALCHANNAG', 7.2235, 4.7472, -4.5034),
    @(12, 'ABREL', 7.2229, 8.0005, 7.553),
    @(13, 'RPOWER', 7.1658, 3.5619, 5.523),
    @(14, 'M&MFIN', 6.6689, 7.295, 16.2005),
    @(15, 'SAIL', 6.5451, 8.7672, 4.707),
    @(16, 'UTKARSHBNK', 6.3655, -4.5601, -1.2393),
    @(17, 'FISCHER', 6.1423, 11.1452, 4.3063),
    @(18, 'SANDUMA', 6.1145, 3.6263, 32.1765),
    @(19, 'ADANIENSOL', 6.003, 3.4154, 11.9518),
    @(20, 'POKARNA', 5.9468, -1.4822, 18.7411),
    @(21, 'JISLJALEQS', 5.8967, 5.1909, -0.8944),
    @(22, 'VAIBHAVGBL', 5.5102, 6.0361, 12.3618),
    @(23, 'CELLO', 5.491, 4.3285, 14.2361),
    @(24, 'GRAPHITE', 5.4633, 11.8789, 12.0803),
    @(25, 'ABDL', 5.3604, 4.2613, 26.9844),
    @(26, 'IOC', 5.2614, 8.1665, 8.5854),
    @(27, 'EPACKPEB', 5.1713, -1.3293, 'N/A'),
    @(28, 'MEGASOFT', 4.9974, 15.7588, 33.5271),
    @(29, 'PROZONER', 4.9921, 15.7468, 36.095),
    @(30, 'STALLION', 4.9914, -5.2229, 21.4391),
    @(31, 'INDOTHAI', 4.9883, 4.7163, 43.9974),
    @(32, 'ATGL', 4.8953, 4.6342, 4.174),
    @(33, 'HITECHGEAR', 4.8651, 2.1287, 10.9905),
    @(34, 'SURYAROSNI', 4.8229, 11.216, 2.8641),
    @(35, 'TMB', 4.7936, 8.5416, 15.7957),
    @(36, 'BAJAJINDEF', 4.7607, 3.6788, 10.7099),
    @(37, 'JKIL', 4.7236, 3.5261, 2.3314),
    @(38, 'GMBREW', 4.5195, 0.0633, 80.0968),
    @(39, 'GENUSPOWER', 4.5101, 2.8253, -0.2258),
    @(40, 'CMSINFO', 4.4508, 3.2219, 3.4293),
    @(41, 'DATAMATICS', 4.4191, 6.8227, 15.1987),
    @(42, 'PDSL', 4.3663, 2.3948, 8.188),
    @(43, 'PROSTARM', 4.3055, 1.4125, -7.5541),
    @(44, 'SAMBHV', 4.2984, 2.7852, 5.3322),
    @(45, 'AXISCADES', 4.2729, 6.7336, -3.2434),
    @(46, 'SGMART', 4.1749, 8.1564, 2.4411),
    @(47, 'LLOYDSENT', 4.0976, 1.3791, 10.7372),
    @(48, 'STLTECH', 4.0809, 1.0909, 7.1775),
    @(49, 'RAJRATAN', 4.0403, 1.5753, 27.754),
    @(50, 'SUNFLAG', 4.0349, 4.371, 4.6693),
    @(51, 'GPPL', 3.9751, 2.9712, 4.6067),
    @(52, 'SRM', 3.959, 3.65, 4.5629),
    @(53, 'GPIL', 3.9558, 6.1219, 14.2198),
    @(54, 'ORIENTTECH', 3.8823, 0.5783, 32.7491),
    @(55, 'RHIM', 3.878, 3.4502, 5.4094),
    @(56, 'HCC', 3.816, 2.717, 7.4566),
    @(57, 'VINCOFE', 3.8087, 10.6775, 9.0515),
    @(58, 'TCI', 3.803, 3.7068, 4.2072),
    @(59, 'NBCC', 3.7531, 2.4711, 6.8827),
    @(60, 'ICRA', 3.7508, 4.4313, 2.8355),
    @(61, 'REDTAPE', 3.7477, 1.6645, -5.14),
    @(62, 'MRPL', 3.7279, 9.146, 19.4367),
    @(63, 'RECLTD', 3.7268, 2.7107, 2.6418),
    @(64, 'MSTCLTD', 3.6113, 3.5433, 15.932),
    @(65, 'STAR', 3.6056, 3.5357, 2.7724),
    @(66, 'ASHAPURMIN', 3.5948, 6.246, 2.0353),
    @(67, 'MAITHANALL', 3.579, 2.6205, 1.8923),
    @(68, 'GAIL', 3.5582, 2.0937, 4.8329),
    @(69, 'DCMSHRIRAM', 3.5476, 10.1086, 17.4754),
    @(70, 'INDORAMA', 3.5365, 2.7019, 13.8319),
    @(71, 'BLACKBUCK', 3.536, 2.154, 8.1153),
    @(72, 'PRAKASH', 3.535, 4.4385, 1.1861),
    @(73, 'SUZLON', 3.5219, 8.1382, 5.7221),
    @(74, 'PVRINOX', 3.5208, 5.6072, 14.0558),
    @(75, 'SHK', 3.4745, 2.2297, -2.0836),
    @(76, 'AVANTEL', 3.4475, -0.0167, 2.9194)
)

foreach ($row in $gainersData) {
    $r = $row[0]
    $gainers.Cells.Item($r, 2).Value = $row[1]
    $gainers.Cells.Item($r, 3).Value = $row[2]
    $gainers.Cells.Item($r, 4).Value = $row[3]
    $gainers.Cells.Item($r, 5).Value = $row[4]
}

# "Top Losers" sheet: CRAMC (row 12) Weekly value update
$losers.Cells.Item(12, 4).Value = 5.9616
